$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = "Sell value actual"
$ws.Range("D2").Value = 15
$ws.Range("D3").Value = 5
$ws.Range("D4").Value = 40

$ws.Columns.Item(4).AutoFit() | Out-Null

$ws.Range("D4").Select() | Out-Null
